$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B10").Value = "24.8650 seconds"
$ws.Range("C10").Value = 2778041
$ws.Range("D10").Value = 1389021
$ws.Range("F10").Value = 262144
$ws.Range("G10").Value = 11156027
$ws.Range("J10").Value = 42.5569

$ws.Columns.Item(4).ColumnWidth = 7.14

$ws.Range("B10").Select()
